# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with the Q1-2022 fund holdings table.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet, placed right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$ws.Name = "2022-Q1"

# Bring over the header/body cell styling (bold centered header with border
# on row 1, bordered index column in column A) from the "2021-Q4" sheet,
# which has the identical 8-column layout.
$q4.Range("B1:H7").Copy($ws.Range("B1:H7"))
$q4.Range("A2:A7").Copy($ws.Range("A2:A7"))

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$rows = @(
    @("162006", "长城久富核心成长混合(LOF)", "19.40", "79.43", "4.47", "0.8672", 6),
    @("008381", "前海开源新兴产业混合", "6.18", "93.63", "7.26", "0.4487", 6),
    @("006976", "鹏华核心优势混合", "2.40", "91.56", "4.54", "0.1090", 5),
    @("000976", "长城新兴产业灵活配置混合", "1.02", "77.95", "4.38", "0.0447", 6),
    @("006923", "前海开源沪港深非周期性行业股票A", "0.54", "93.77", "6.57", "0.0355", 3),
    @("006924", "前海开源沪港深非周期性行业股票C", "0.22", "93.77", "6.57", "0.0145", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    # Fund codes and the decimal-formatted figures are stored as plain text
    # in the source data (e.g. "008381", "19.40"), so force text entry with
    # a leading quote-prefix instead of letting them coerce to numbers.
    $ws.Cells.Item($r, 2).Value = "'" + $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = "'" + $row[2]
    $ws.Cells.Item($r, 5).Value = "'" + $row[3]
    $ws.Cells.Item($r, 6).Value = "'" + $row[4]
    $ws.Cells.Item($r, 7).Value = "'" + $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row on the "总计" sheet
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Row 4 doesn't exist yet - seed it (and its styling) from row 3, which
# already carries the right "index column" border/format on A3.
$zj.Range("A3:D3").Copy($zj.Range("A4:D4"))

# Push the existing data rows down by one (bottom-up to avoid clobbering).
# Column A is just a 0-based row counter (0,1,2,...), so it is left as-is -
# only the data columns B:D need to move.
$zj.Range("B3:D3").Copy($zj.Range("B4:D4"))
$zj.Range("B2:D2").Copy($zj.Range("B3:D3"))
$zj.Cells.Item(4, 1).Value = 2

$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 6
$zj.Cells.Item(2, 4).Value = 1.52

# Restore the originally active sheet (adding a sheet shifts focus to it).
$wb.Worksheets.Item("2021-Q3").Select()
